$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is published, so it is inserted as the first data row
# (row 14) of this subset, pushing all the subsequent rows down by one and
# growing the used range from A1:R39 to A1:R40.
$ws.Rows(14).Insert()

# Fill in the freshly inserted row with the new record's data.
$ws.Cells.Item(14, 1).Value = 3
$ws.Cells.Item(14, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(14, 3).Value = "Coquimbo"
$ws.Cells.Item(14, 4).Value = 44525
$ws.Cells.Item(14, 5).Value = 5
$ws.Cells.Item(14, 6).Value = 100112022
$ws.Cells.Item(14, 7).Value = "Arveja Verde"
$ws.Cells.Item(14, 8).Value = "Perfection"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 73
$ws.Cells.Item(14, 11).Value = 16000
$ws.Cells.Item(14, 12).Value = 17000
$ws.Cells.Item(14, 13).Value = 16479
$ws.Cells.Item(14, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(14, 15).Value = "Provincia de Talca"
$ws.Cells.Item(14, 16).Value = 659
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"
